$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.524740333333333
$ws.Range("H2").Value = 4.574221
$ws.Range("I2").Value = 0.2062237893390968
$ws.Range("J2").Value = 0.2062237893390969
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06024766666666666
$ws.Range("N2").Value = 0.180743
$ws.Range("O2").Value = 0.03337856003285722
$ws.Range("P2").Value = 0.03337856003285723
$ws.Range("Q2").Value = 0.09186204735588888
$ws.Range("R2").Value = 0.8267584262029999
$ws.Range("S2").Value = 0.006883453132658344
$ws.Range("T2").Value = 0.006883453132658347

$ws.Range("G3").Value = 1.524740333333333
$ws.Range("H3").Value = 4.574221
$ws.Range("I3").Value = 0.2062237893390968
$ws.Range("J3").Value = 0.2062237893390969
$ws.Range("N3").Value = 5.09203
$ws.Range("O3").Value = 0.9403663159519869
$ws.Range("P3").Value = 0.9403663159519869
$ws.Range("Q3").Value = 2.588007839847778
$ws.Range("R3").Value = 23.29207055863
$ws.Range("S3").Value = 0.1939259050424651
$ws.Range("T3").Value = 0.1939259050424652

$ws.Range("G4").Value = 1.524740333333333
$ws.Range("H4").Value = 4.574221
$ws.Range("I4").Value = 0.2062237893390968
$ws.Range("J4").Value = 0.2062237893390969
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.04738999999999999
$ws.Range("N4").Value = 0.14217
$ws.Range("O4").Value = 0.02625512401515583
$ws.Range("P4").Value = 0.02625512401515584
$ws.Range("Q4").Value = 0.07225744439666665
$ws.Range("R4").Value = 0.6503169995699999
$ws.Range("S4").Value = 0.005414431163973358
$ws.Range("T4").Value = 0.005414431163973361

$ws.Range("I5").Value = 0.4308548451232278
$ws.Range("J5").Value = 0.4308548451232279
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06024766666666666
$ws.Range("N5").Value = 0.180743
$ws.Range("O5").Value = 0.03337856003285722
$ws.Range("P5").Value = 0.03337856003285723
$ws.Range("Q5").Value = 0.191923581237
$ws.Range("R5").Value = 1.727312231133
$ws.Range("S5").Value = 0.01438131431339306
$ws.Range("T5").Value = 0.01438131431339307

$ws.Range("I6").Value = 0.4308548451232278
$ws.Range("J6").Value = 0.4308548451232279
$ws.Range("N6").Value = 5.09203
$ws.Range("O6").Value = 0.9403663159519869
$ws.Range("P6").Value = 0.9403663159519869
$ws.Range("Q6").Value = 5.40701788377
$ws.Range("R6").Value = 48.66316095393
$ws.Range("S6").Value = 0.4051613834185936
$ws.Range("T6").Value = 0.4051613834185938

$ws.Range("I7").Value = 0.4308548451232278
$ws.Range("J7").Value = 0.4308548451232279
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.04738999999999999
$ws.Range("N7").Value = 0.14217
$ws.Range("O7").Value = 0.02625512401515583
$ws.Range("P7").Value = 0.02625512401515584
$ws.Range("Q7").Value = 0.15096449403
$ws.Range("R7").Value = 1.35868044627
$ws.Range("S7").Value = 0.01131214739124111
$ws.Range("T7").Value = 0.01131214739124111

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2032796666666667
$ws.Range("H8").Value = 0.609839
$ws.Range("I8").Value = 0.02749392945088694
$ws.Range("J8").Value = 0.02749392945088694
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.06024766666666666
$ws.Range("N8").Value = 0.180743
$ws.Range("O8").Value = 0.03337856003285722
$ws.Range("P8").Value = 0.03337856003285723
$ws.Range("Q8").Value = 0.01224712559744444
$ws.Range("R8").Value = 0.110224130377
$ws.Range("S8").Value = 0.0009177077747155707
$ws.Range("T8").Value = 0.0009177077747155711

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2032796666666667
$ws.Range("H9").Value = 0.609839
$ws.Range("I9").Value = 0.02749392945088694
$ws.Range("J9").Value = 0.02749392945088694
$ws.Range("N9").Value = 5.09203
$ws.Range("O9").Value = 0.9403663159519869
$ws.Range("P9").Value = 0.9403663159519869
$ws.Range("Q9").Value = 0.3450353870188889
$ws.Range("R9").Value = 3.10531848317
$ws.Range("S9").Value = 0.02585436514877438
$ws.Range("T9").Value = 0.02585436514877439

$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.2032796666666667
$ws.Range("H10").Value = 0.609839
$ws.Range("I10").Value = 0.02749392945088694
$ws.Range("J10").Value = 0.02749392945088694
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.04738999999999999
$ws.Range("N10").Value = 0.14217
$ws.Range("O10").Value = 0.02625512401515583
$ws.Range("P10").Value = 0.02625512401515584
$ws.Range("Q10").Value = 0.009633423403333332
$ws.Range("R10").Value = 0.08670081063
$ws.Range("S10").Value = 0.0007218565273969818
$ws.Range("T10").Value = 0.0007218565273969821

$ws.Range("G11").Value = 2.061212666666667
$ws.Range("H11").Value = 6.183638
$ws.Range("I11").Value = 0.2787826080683977
$ws.Range("J11").Value = 0.2787826080683978
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.06024766666666666
$ws.Range("N11").Value = 0.180743
$ws.Range("O11").Value = 0.03337856003285722
$ws.Range("P11").Value = 0.03337856003285723
$ws.Range("Q11").Value = 0.1241832536704444
$ws.Range("R11").Value = 1.117649283034
$ws.Range("S11").Value = 0.009305362019527516
$ws.Range("T11").Value = 0.009305362019527523

$ws.Range("G12").Value = 2.061212666666667
$ws.Range("H12").Value = 6.183638
$ws.Range("I12").Value = 0.2787826080683977
$ws.Range("J12").Value = 0.2787826080683978
$ws.Range("N12").Value = 5.09203
$ws.Range("O12").Value = 0.9403663159519869
$ws.Range("P12").Value = 0.9403663159519869
$ws.Range("Q12").Value = 3.498585578348889
$ws.Range("R12").Value = 31.48727020514
$ws.Range("S12").Value = 0.2621577741007658
$ws.Range("T12").Value = 0.2621577741007659

$ws.Range("G13").Value = 2.061212666666667
$ws.Range("H13").Value = 6.183638
$ws.Range("I13").Value = 0.2787826080683977
$ws.Range("J13").Value = 0.2787826080683978
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.04738999999999999
$ws.Range("N13").Value = 0.14217
$ws.Range("O13").Value = 0.02625512401515583
$ws.Range("P13").Value = 0.02625512401515584
$ws.Range("Q13").Value = 0.09768086827333332
$ws.Range("R13").Value = 0.87912781446
$ws.Range("S13").Value = 0.007319471948104363
$ws.Range("T13").Value = 0.007319471948104368

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4188103333333333
$ws.Range("H14").Value = 1.256431
$ws.Range("I14").Value = 0.05664482801839063
$ws.Range("J14").Value = 0.05664482801839064
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.06024766666666666
$ws.Range("N14").Value = 0.180743
$ws.Range("O14").Value = 0.03337856003285722
$ws.Range("P14").Value = 0.03337856003285723
$ws.Range("Q14").Value = 0.02523234535922222
$ws.Range("R14").Value = 0.227091108233
$ws.Range("S14").Value = 0.001890722792562724
$ws.Range("T14").Value = 0.001890722792562725

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4188103333333333
$ws.Range("H15").Value = 1.256431
$ws.Range("I15").Value = 0.05664482801839063
$ws.Range("J15").Value = 0.05664482801839064
$ws.Range("N15").Value = 5.09203
$ws.Range("O15").Value = 0.9403663159519869
$ws.Range("P15").Value = 0.9403663159519869
$ws.Range("Q15").Value = 0.7108649272144445
$ws.Range("R15").Value = 6.397784344930001
$ws.Range("S15").Value = 0.05326688824138789
$ws.Range("T15").Value = 0.0532668882413879

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4188103333333333
$ws.Range("H16").Value = 1.256431
$ws.Range("I16").Value = 0.05664482801839063
$ws.Range("J16").Value = 0.05664482801839064
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.04738999999999999
$ws.Range("N16").Value = 0.14217
$ws.Range("O16").Value = 0.02625512401515583
$ws.Range("P16").Value = 0.02625512401515584
$ws.Range("Q16").Value = 0.01984742169666667
$ws.Range("R16").Value = 0.17862679527
$ws.Range("S16").Value = 0.00148721698444002
$ws.Range("T16").Value = 0.00148721698444002
